# PUBLONS.xlsx edit: add PUBLONS002-005 test cases + new PUBLONS005 sheet
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Test Cases" sheet: fill rows 3-5 (PUBLONS002, PUBLONS003, PUBLONS004)
#    Columns: A=TCID, B=JIRA ID, C=Description, D=Runmode
#    Values are written in this specific order to reproduce the shared
#    string table order of the original author's edit.
# ---------------------------------------------------------------------
$tc = $wb.Worksheets.Item("Test Cases")

$tc.Range("A3").Value = "PUBLONS002"
$tc.Range("A4").Value = "PUBLONS003"
$tc.Range("B3").Value = "OPQA-5765"
$tc.Range("C3").Value = "Verify that email address field should be in standard email ID format .Email address fields should be mandatory."
$tc.Range("C4").Value = "Verify that error message ""Please enter an email address"" whenever not enter any text in email field"
$tc.Range("B4").Value = "OPQA-5766"
$tc.Range("A5").Value = "PUBLONS004"
$tc.Range("B5").Value = "OPQA-5767"
$tc.Range("C5").Value = "Verify that error message ""Please enter a valid email address"" whenever enter wrong format in email field"

$tc.Range("D3").Value = "Y"
$tc.Range("D4").Value = "Y"
$tc.Range("D5").Value = "Y"

# ---------------------------------------------------------------------
# 2. New worksheet "PUBLONS005", inserted right before "WAT09"
# ---------------------------------------------------------------------
$wat09 = $wb.Worksheets.Item("WAT09")
$ws = $wb.Worksheets.Add($wat09)
$ws.Name = "PUBLONS005"

$ws.Range("A1").Value = "CHARACTER LENGTH"
$ws.Range("B1").Value = "SUFFIX"
$ws.Range("C1").Value = "ERROR TEXT"
$ws.Range("D1").Value = "VALIDITY"
$ws.Range("E1").Value = "Runmode"
$ws.Range("F1").Value = "PASS"

# Mailto hyperlinks on the three numeric "CHARACTER LENGTH" cells (added
# before the numeric values so the hyperlink display text doesn't clobber
# the cell's numeric content).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:amneetsingh72@gmail.com", [Type]::Missing, [Type]::Missing, "amneetsingh72@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:amneetsinghasr@gmail.com", [Type]::Missing, [Type]::Missing, "amneetsinghasr@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:amneetsingh100@gmail.com", [Type]::Missing, [Type]::Missing, "amneetsingh100@gmail.com")

$ws.Range("A2").Value = 246
$ws.Range("B2").Value = "@abc.com"
$ws.Range("D2").Value = "YES"
$ws.Range("E2").Value = "Y"
$ws.Range("F2").Value = "SKIP"

$ws.Range("A3").Value = 247
$ws.Range("B3").Value = "@abc.com"
$ws.Range("D3").Value = "YES"
$ws.Range("E3").Value = "Y"
$ws.Range("F3").Value = "SKIP"

$ws.Range("A4").Value = 248
$ws.Range("B4").Value = "@abc.com"
$ws.Range("C4").Value = "Please enter no more than 255 characters."
$ws.Range("D4").Value = "NO"
$ws.Range("E4").Value = "Y"
$ws.Range("F4").Value = "SKIP"

# Borders around the used range, matching the bordered-table look of
# the rest of the workbook.
$ws.Range("A1:F4").Borders.LineStyle = 1

# Column widths (best-fit approximations)
$ws.Columns.Item(1).ColumnWidth = 18.17
$ws.Columns.Item(2).ColumnWidth = 9.33
$ws.Columns.Item(4).ColumnWidth = 8
$ws.Columns.Item(6).ColumnWidth = 4.67

$ws.Range("A4").Select()

# ---------------------------------------------------------------------
# 3. "Test Cases" sheet: fill row 6 (PUBLONS005) - added after the new
#    sheet so the shared strings land at the end of the table, matching
#    the source edit.
# ---------------------------------------------------------------------
$tc.Range("A6").Value = "PUBLONS005"
$tc.Range("C6").Value = "Verify that error message ""Email address is too long."" whenever enter more than 255 characters in email field&&Verify that email address field should be in standard email ID format .Email address fields should be mandatory."
$tc.Range("B6").Value = "OPQA-5768&&OPQA-5765"
$tc.Range("D6").Value = "Y"

$tc.Rows.Item(6).RowHeight = 30

$tc.Range("B6").Select()
